# Update the "具体时间范围" (E) and "想去人数" (F) columns on the
# 展览 and 全部类型 sheets: normalize "HH:MM - MM.DD HH:MM" strings to
# "HH:MM-MM.DD HH:MM" (remove spaces around the dash) and bump a handful
# of the attendee counts.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# row => (new E text, new F value-or-$null-if-unchanged)
$rowUpdates = @{
    2  = @{ E = "2024.03.16 13:00-03.16 17:30"; F = 160 }
    3  = @{ E = "2024.03.17 10:00-03.17 17:00"; F = 1739 }
    4  = @{ E = "2024.03.23 10:00-03.23 17:30"; F = 796 }
    5  = @{ E = "2024.03.30 10:00-03.30 17:00"; F = 2 }
    6  = @{ E = "2024.03.30 10:00-03.31 17:00"; F = $null }
    7  = @{ E = "2024.03.31 10:00-03.31 17:00"; F = $null }
    8  = @{ E = "2024.04.04 10:00-04.05 17:00"; F = 12067 }
    9  = @{ E = "2024.04.04 10:00-04.04 17:00"; F = $null }
    10 = @{ E = "2024.04.06 10:00-04.06 16:00"; F = 100 }
    11 = @{ E = "2024.04.13 10:00-04.14 17:00"; F = $null }
    12 = @{ E = "2024.04.13 13:00-04.13 20:00"; F = 419 }
    13 = @{ E = "2024.04.20 10:00-04.20 17:00"; F = $null }
    14 = @{ E = "2024.04.21 10:00-04.21 21:00"; F = $null }
    15 = @{ E = "2024.05.01 10:00-05.03 17:00"; F = 13509 }
    16 = @{ E = "2024.05.01 10:00-05.02 17:00"; F = 13546 }
    17 = @{ E = "2024.05.02 14:00-05.02 16:00"; F = $null }
    18 = @{ E = "2024.05.02 14:00-05.02 16:00"; F = $null }
    19 = @{ E = "2024.05.02 14:00-05.02 16:00"; F = $null }
    20 = @{ E = "2024.05.02 14:00-05.02 16:00"; F = $null }
    21 = @{ E = "2024.05.03 09:00-05.04 17:30"; F = 994 }
    22 = @{ E = "2024.05.03 14:00-05.03 16:00"; F = $null }
    23 = @{ E = "2024.05.03 14:00-05.03 16:00"; F = 50 }
    24 = @{ E = "2024.05.04 09:00-05.05 17:00"; F = 1975 }
    25 = @{ E = "2024.06.08 10:00-06.09 17:00"; F = 181 }
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $rowUpdates.Keys) {
        $update = $rowUpdates[$row]

        $ws.Cells.Item($row, 5).Value = $update.E

        if ($null -ne $update.F) {
            $ws.Cells.Item($row, 6).Value = $update.F
        }
    }
}
